$wb = $excel.ActiveWorkbook

# ---- Sheet "Login": add new test case row, move the "ON" flag down ----
$ws1 = $wb.Worksheets.Item("Login")

# Row 7 currently holds ("ON","SF_Sign_In_TC006",1) with A7 unstyled.
# Give A7 the same border style as the rows above it, then move its "ON"
# text down onto the new row 8.
$ws1.Cells.Item(2,1).Copy()
$ws1.Cells.Item(7,1).PasteSpecial(-4122)
$ws1.Cells.Item(7,1).ClearContents()

# Highlight the test-case-name cell (B7) with a themed fill and a full
# thin border (it used to be left/right-only).
$ws1.Cells.Item(7,2).Interior.Color = 12444887
$ws1.Cells.Item(7,2).Borders.Color = 0

# New row 8: ON | SF_Shopping_Sign_In_TC007 | 1 -- clone formatting from
# the row above (A/C from row 2, B from the now-restyled row 7) so the
# new row matches the existing look.
$ws1.Cells.Item(2,1).Copy()
$ws1.Cells.Item(8,1).PasteSpecial(-4122)
$ws1.Cells.Item(2,3).Copy()
$ws1.Cells.Item(8,3).PasteSpecial(-4122)
$ws1.Cells.Item(7,2).Copy()
$ws1.Cells.Item(8,2).PasteSpecial(-4122)

$ws1.Cells.Item(8,1).Value = "ON"
$ws1.Cells.Item(8,2).Value = "SF_Shopping_Sign_In_TC007"
$ws1.Cells.Item(8,3).Value = 1

# Column widths
$ws1.Columns.Item(1).ColumnWidth = 3.584
$ws1.Columns.Item(2).ColumnWidth = 44.75
$ws1.Columns.Item(3).ColumnWidth = 6.584

# View: scroll down a bit and land the selection on the newly added cell
$ws1.Activate()
$ws1.Range("B8").Select()

# ---- Sheet "Jobs": column B width tweak ----
$ws2 = $wb.Worksheets.Item("Jobs")
$ws2.Columns.Item(2).ColumnWidth = 51.584
